# Add IEC 5 minute weatherfiles - rows 74-77 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C (IEC method names) first, top to bottom,
#     so the shared-string table gets IEC1B, IEC2Mono, IEC2A, IEC2B in that order ---
$ws.Range("B74").Value = "IEC1B"
$ws.Range("C74").Value = "'IEC1B"

$ws.Range("B75").Value = "'IEC2Mono"
$ws.Range("C75").Value = "'IEC2Mono"

$ws.Range("B76").Value = "'IEC2A"
$ws.Range("C76").Value = "'IEC2A"

$ws.Range("B77").Value = "'IEC2B"
$ws.Range("C77").Value = "'IEC2B"

# --- Column D (description) entered bottom to top, matching the
#     shared-string insertion order seen in the target file ---
$ws.Range("D77").Value = "5-minute data - refmod"
$ws.Range("D76").Value = "5-minute data - POA+BOA"
$ws.Range("D75").Value = "5-minute data POA"
$ws.Range("D74").Value = "5-minute data GHI+DHI"

# --- Remaining columns, row by row ---
$ws.Range("A74").Value = 1
$ws.Range("E74").Value = "SRRL_DNI"
$ws.Range("F74").Value = "SRRL_DHI"
$ws.Range("G74").Value = "SRRL_GHI"
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("O74").Value = 2
$ws.Range("P74").Value = 2
$ws.Range("Q74").Value = "All tree"
$ws.Range("S74").Value = "Yes"
$ws.Range("T74").Value = "Yes"
$ws.Range("U74").Value = "Yes"
$ws.Range("V74").Value = "Yes"

$ws.Range("A75").Value = 1
$ws.Range("E75").Value = "SRRL_DNI"
$ws.Range("F75").Value = "SRRL_DHI"
$ws.Range("G75").Value = "SRRL_GHI"
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("O75").Value = 4
$ws.Range("P75").Value = 2
$ws.Range("Q75").Value = "All tree"
$ws.Range("S75").Value = "Yes"
$ws.Range("T75").Value = "Yes"
$ws.Range("U75").Value = "Yes"
$ws.Range("V75").Value = "Yes"

$ws.Range("A76").Value = 1
$ws.Range("E76").Value = "SRRL_DNI"
$ws.Range("F76").Value = "SRRL_DHI"
$ws.Range("G76").Value = "SRRL_GHI"
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 2
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("O76").Value = 4
$ws.Range("P76").Value = 2
$ws.Range("Q76").Value = "All tree"
$ws.Range("S76").Value = "Yes"
$ws.Range("T76").Value = "Yes"
$ws.Range("U76").Value = "Yes"
$ws.Range("V76").Value = "Yes"

$ws.Range("A77").Value = 1
$ws.Range("E77").Value = "SRRL_DNI"
$ws.Range("F77").Value = "SRRL_DHI"
$ws.Range("G77").Value = "SRRL_GHI"
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 2
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("O77").Value = 4
$ws.Range("P77").Value = 2
$ws.Range("Q77").Value = "All tree"
$ws.Range("S77").Value = "Yes"
$ws.Range("T77").Value = "Yes"
$ws.Range("U77").Value = "Yes"
$ws.Range("V77").Value = "Yes"

# Apply the cell formatting/styles matching existing rows in this table
$ws.Range("G24").Copy()
$ws.Range("G74:G77").PasteSpecial(-4122)
$ws.Range("H24:J24").Copy()
$ws.Range("H74:J77").PasteSpecial(-4122)

# Leave the selection where the author would have ended up after typing
# the new rows (just below the last entry), and scroll the view down.
$excel.ActiveWindow.ScrollRow = 65
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D79").Select() | Out-Null
